# Add new match-record rows to the "Partidos" sheet (rows 456-470),
# matching the data appended in the source workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

# Each entry: row, jugador, equipo, posicion, goles, autogoles, arquero(bool),
#             goles_recibidos, tarjetas_amarillas, tarjetas_rojas, asistencias, penales_atajados
$newRows = @(
    @(456, "Invitado",              "Azul",     "Arquero",       0, 0, $true,  5, 0, 0, 0, 0),
    @(457, "Gember Marin Sarria",   "Amarillo", "Arquero",       0, 0, $true,  4, 0, 0, 0, 0),
    @(458, "Fede",                  "Amarillo", "Arquero",       0, 0, $true,  1, 0, 0, 0, 0),
    @(459, "Juan David Espinal",    "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @(460, "Juan Felipe Gutierrez", "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 1, 0),
    @(461, "Cesar Augusto Estrada", "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @(462, "Armando Murillo",       "Azul",     "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @(463, "Arnul David Narvaez",   "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @(464, "Fabian Grajales",       "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(465, "Juan Carlos Otero",     "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(466, "Alexander Uribe",       "Amarillo", "Mediocampista", 2, 0, $false, 0, 0, 0, 0, 0),
    @(467, "Andres Tangarife",      "Amarillo", "Delantero",     3, 0, $false, 0, 0, 0, 0, 0),
    @(468, "Andres Jurado",         "Amarillo", "Delantero",     0, 0, $false, 0, 0, 0, 3, 0),
    @(469, "Bryan Andres Burgos",   "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @(470, "Armando Vieras",        "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0)
)

$fechaPartido = 45934

foreach ($r in $newRows) {
    $rowNum = $r[0]

    $ws.Cells.Item($rowNum, 1).Value = $fechaPartido
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
    $ws.Cells.Item($rowNum, 6).Value = $r[5]
    $ws.Cells.Item($rowNum, 7).Value = $r[6]
    $ws.Cells.Item($rowNum, 8).Value = $r[7]
    $ws.Cells.Item($rowNum, 9).Value = $r[8]
    $ws.Cells.Item($rowNum, 10).Value = $r[9]
    $ws.Cells.Item($rowNum, 11).Value = $r[10]
    $ws.Cells.Item($rowNum, 12).Value = $r[11]
}

# Leave the selection where the author ended up after entering the data.
$ws.Activate()
$ws.Range("B471").Select()
